$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (was "Исходный ключ" -> now "Пин-код"); B3:M3 values unchanged but highlighted yellow ---
$ws.Range("A3").Value2 = "Пин-код"

# --- Row 6 becomes "Контрольное число текущего XOR ключа" header row (moved from old row 8) ---
$ws.Range("A6").Value2 = "Контрольное число текущего XOR ключа"
$ws.Range("B6").Formula = "=MOD(SUM(B5:M5),12)"
$ws.Range("C6:M6").ClearContents()
$ws.Range("B6:M6").Merge()

# --- Row 7 becomes "Ключ сдвига" row (moved from old row 6), formulas updated to reference $B$6 with IF/LEFT guard
$ws.Range("A7").Value2 = "Ключ сдвига"
$ws.Range("B7").Formula = "=IF(B5>1000,LEFT(B5,3)+B2,B5)+`$B`$6"
$ws.Range("C7").Formula = "=IF(C5>1000,LEFT(C5,3)+C2,C5)-`$B`$6"
$ws.Range("D7").Formula = "=IF(D5>1000,LEFT(D5,3)+D2,D5)+`$B`$6"
$ws.Range("E7").Formula = "=IF(E5>1000,LEFT(E5,3)+E2,E5)-`$B`$6"
$ws.Range("F7").Formula = "=IF(F5>1000,LEFT(F5,3)+F2,F5)+`$B`$6"
$ws.Range("G7").Formula = "=IF(G5>1000,LEFT(G5,3)+G2,G5)-`$B`$6"
$ws.Range("H7").Formula = "=IF(H5>1000,LEFT(H5,3)+H2,H5)+`$B`$6"
$ws.Range("I7").Formula = "=IF(I5>1000,LEFT(I5,3)+I2,I5)-`$B`$6"
$ws.Range("J7").Formula = "=IF(J5>1000,LEFT(J5,3)+J2,J5)+`$B`$6"
$ws.Range("K7").Formula = "=IF(K5>1000,LEFT(K5,3)+K2,K5)-`$B`$6"
$ws.Range("L7").Formula = "=IF(L5>1000,LEFT(L5,3)+L2,L5)+`$B`$6"
$ws.Range("M7").Formula = "=IF(M5>1000,LEFT(M5,3)+M2,M5)-`$B`$6"

# --- Row 8: old filler row becomes mostly empty (A8 kept bordered/empty; old content cleared) ---
$ws.Range("A8").ClearContents()
$ws.Range("B8:M8").ClearContents()

# --- Row 9: fully cleared (old "Пин-код"/12345 row removed) ---
$ws.Range("A9:N9").Clear()

Write-Host "done"
